$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'19.969.48"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "'1.414.55"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'276.11"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "'0.3690"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "'0.3108"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "'39.83"
$ws.Range("E9").Value = "  -3.77%  "
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").Value = "'0.06508"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'5.468"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'17.60"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'6.180"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'1.414.82"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'0.05679"
$ws.Range("E18").Value = "  -6.22%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "'70.92"
$ws.Range("E20").Value = "  -8.88%  "
$ws.Range("D21").Value = "'5.588"
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'10.98"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'2.237"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "'19.971.90"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").Value = "'2.275"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "'132.97"
$ws.Range("E27").Value = "  -7.19%  "
$ws.Range("D28").Value = "'17.20"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'1.574.86"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").Value = "'109.73"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'3.919"
$ws.Range("E31").Value = "  +5.02%  "
$ws.Range("D32").Value = "'5.193"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "'0.8098"
$ws.Range("E33").Value = "  -11.57%  "
$ws.Range("D34").Value = "'0.07768"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'1.468"
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'4.878"
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05817"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'8.044"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").Value = "'1.000"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'0.02045"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'10.44"
$ws.Range("E41").Value = "  -5.27%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1883"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.101"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.5299"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.34"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").Value = "'3.535"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "'116.66"
$ws.Range("E47").Value = "  +5.12%  "
$ws.Range("D48").Value = "'0.5180"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "'1.766"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "'1.032"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.74%  "
